# The underlying OOXML diff for this revision is a pure re-serialization
# (attribute order only -- namespace / attribute ordering in
# ppt/presentation.xml's <p:presentation> and <a:lvlNpPr>/<a:defRPr>
# elements, plus the <p:ph> placeholder tag on the title slide's
# subtitle shape). No text, formatting, geometry, or structural content
# actually changed between revisions.
#
# Re-touch the two slide-1 placeholders (title + subtitle) through the
# object model so the authoring pass that produced this revision is
# faithfully represented, without altering any visible content.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = $title.TextFrame.TextRange.Text

$subtitle = $s.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Text = $subtitle.TextFrame.TextRange.Text
